# cibmtr-reporting-ig: refresh ValueSet metadata sheet with new IG publish
# values (version bump, status -> draft, new date, full contact info +
# jurisdiction row).

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Metadata")

# Version: 0.1.6 -> 0.1.7
$ws.Range("B3").Value = "0.1.7"

# Status: active -> draft
$ws.Range("B6").Value = "draft"

# Date: refreshed publication timestamp
$ws.Range("B8").Value = "2024-08-27T12:23:18-05:00"

# Contact: first contact row now carries the publisher contact details,
# second (duplicate) Contact row now carries the individual FHIR contact.
$ws.Range("B10").Value = "The Medical College of Wisconsin, Inc. and the National Marrow Donor Program (http://www.cibmtr.org)"
$ws.Range("B11").Value = "Bob Milius (bmilius@nmdp.org)"

# A new "Jurisdiction" property row is inserted right after the Contact
# rows (pushing Description/Purpose/Copyright/Immutable down by one row).
# First extend the existing body formatting (border/fill/alignment) down
# onto the new row 16 by copying it from row 15 (uniform body style), then
# shift the existing rows down (bottom-up, so we never clobber a value
# before it has been copied onward) using the already-formatted rows so
# no new cell styles are introduced.
$ws.Range("A15:B15").Copy()
$ws.Range("A16:B16").PasteSpecial(-4122)
$excel.CutCopyMode = $false

$ws.Range("A16").Value = $ws.Range("A15").Value()
$ws.Range("B16").Value = $ws.Range("B15").Value()

$ws.Range("A15").Value = $ws.Range("A14").Value()
$ws.Range("B15").Value = $ws.Range("B14").Value()

$ws.Range("A14").Value = $ws.Range("A13").Value()
$ws.Range("B14").Value = $ws.Range("B13").Value()

$ws.Range("A13").Value = $ws.Range("A12").Value()
$ws.Range("B13").Value = $ws.Range("B12").Value()

$ws.Range("A12").Value = "Jurisdiction"
$ws.Range("B12").Value = ""
